# Replace the placeholder numeric examples in each lookup sheet with
# more meaningful example values (commit: "change to better examples").

$wb = $excel.ActiveWorkbook

# --- system sheet: LNC/Linac, BSTR/Booster -------------------------------
$ws = $wb.Worksheets.Item("system")
$ws.Range("A2").Value = "LNC"
$ws.Range("B2").Value = "Linac"
$ws.Range("A3").Value = "BSTR"
$ws.Range("B3").Value = "Booster"
[void]$ws.Range("C4").Select()

# --- subsystem sheet: MAG/magnet ------------------------------------------
$ws = $wb.Worksheets.Item("subsystem")
$ws.Range("A2").Value = "MAG"
$ws.Range("B2").Value = "magnet"
[void]$ws.Range("A3").Select()

# --- device_type sheet: BEND/dipole bend, QUAD/quadrupole -----------------
$ws = $wb.Worksheets.Item("device_type")
$ws.Range("A2").Value = "BEND"
$ws.Range("B2").Value = "dipole bend"
$ws.Range("A3").Value = "QUAD"
$ws.Range("B3").Value = "quadrupole"
[void]$ws.Range("B4").Select()

# --- signal_type sheet: Amp_Set/amplitude set point, I_Set/current set point
$ws = $wb.Worksheets.Item("signal_type")
$ws.Range("A2").Value = "Amp_Set"
$ws.Range("B2").Value = "amplitude set point"
$ws.Range("A3").Value = "I_Set"
$ws.Range("B3").Value = "current set point"
[void]$ws.Range("B4").Select()
